# Update "想去人数" (F column) counts and mark one ticket as "不可售" (G column)
# on the "展览" and "全部类型" worksheets, per the upstream gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Sheets.Item("展览")
$wsExpo.Cells.Item(2, 6).Value = 125
$wsExpo.Cells.Item(2, 7).Value = "不可售"
$wsExpo.Cells.Item(3, 6).Value = 265
$wsExpo.Cells.Item(4, 6).Value = 149
$wsExpo.Cells.Item(5, 6).Value = 1808
$wsExpo.Cells.Item(6, 6).Value = 1543
$wsExpo.Cells.Item(7, 6).Value = 279
$wsExpo.Cells.Item(8, 6).Value = 67
$wsExpo.Cells.Item(9, 6).Value = 516
$wsExpo.Cells.Item(10, 6).Value = 128

$wsAll = $wb.Sheets.Item("全部类型")
$wsAll.Cells.Item(2, 6).Value = 125
$wsAll.Cells.Item(2, 7).Value = "不可售"
$wsAll.Cells.Item(3, 6).Value = 265
$wsAll.Cells.Item(4, 6).Value = 149
$wsAll.Cells.Item(5, 6).Value = 1808
$wsAll.Cells.Item(6, 6).Value = 1543
$wsAll.Cells.Item(7, 6).Value = 279
$wsAll.Cells.Item(9, 6).Value = 67
$wsAll.Cells.Item(10, 6).Value = 516
$wsAll.Cells.Item(11, 6).Value = 128
